# Update NATMI TPM-derived columns (G:J and M:T) for rows 2-10 on the active sheet
# with newly recomputed TPM-based values, per commit "update scripts wuth new tpm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.12586833333333
$ws.Range("H2").Value = 30.377605
$ws.Range("I2").Value = 0.9311967029481902
$ws.Range("J2").Value = 0.9311967029481902
$ws.Range("M2").Value = 1.492477333333333
$ws.Range("N2").Value = 4.477432
$ws.Range("O2").Value = 0.02769484181536182
$ws.Range("P2").Value = 0.02769484181536182
$ws.Range("Q2").Value = 15.11262896781778
$ws.Range("R2").Value = 136.01366071036
$ws.Range("S2").Value = 0.0257893453871366
$ws.Range("T2").Value = 0.0257893453871366

$ws.Range("G3").Value = 10.12586833333333
$ws.Range("H3").Value = 30.377605
$ws.Range("I3").Value = 0.9311967029481902
$ws.Range("J3").Value = 0.9311967029481902
$ws.Range("M3").Value = 33.85786133333334
$ws.Range("N3").Value = 101.573584
$ws.Range("O3").Value = 0.6282762845978157
$ws.Range("P3").Value = 0.6282762845978156
$ws.Range("Q3").Value = 342.8402459095912
$ws.Range("R3").Value = 3085.562213186321
$ws.Range("S3").Value = 0.5850488047580248
$ws.Range("T3").Value = 0.5850488047580247

$ws.Range("G4").Value = 10.12586833333333
$ws.Range("H4").Value = 30.377605
$ws.Range("I4").Value = 0.9311967029481902
$ws.Range("J4").Value = 0.9311967029481902
$ws.Range("M4").Value = 18.53974466666667
$ws.Range("N4").Value = 55.619234
$ws.Range("O4").Value = 0.3440288735868225
$ws.Range("P4").Value = 0.3440288735868225
$ws.Range("Q4").Value = 187.7310134282856
$ws.Range("R4").Value = 1689.57912085457
$ws.Range("S4").Value = 0.3203585528030289
$ws.Range("T4").Value = 0.3203585528030289

$ws.Range("G5").Value = 0.100996
$ws.Range("H5").Value = 0.302988
$ws.Range("I5").Value = 0.009287810103293732
$ws.Range("J5").Value = 0.009287810103293733
$ws.Range("M5").Value = 1.492477333333333
$ws.Range("N5").Value = 4.477432
$ws.Range("O5").Value = 0.02769484181536182
$ws.Range("P5").Value = 0.02769484181536182
$ws.Range("Q5").Value = 0.1507342407573333
$ws.Range("R5").Value = 1.356608166816
$ws.Range("S5").Value = 0.0002572244316218393
$ws.Range("T5").Value = 0.0002572244316218393

$ws.Range("G6").Value = 0.100996
$ws.Range("H6").Value = 0.302988
$ws.Range("I6").Value = 0.009287810103293732
$ws.Range("J6").Value = 0.009287810103293733
$ws.Range("M6").Value = 33.85786133333334
$ws.Range("N6").Value = 101.573584
$ws.Range("O6").Value = 0.6282762845978157
$ws.Range("P6").Value = 0.6282762845978156
$ws.Range("Q6").Value = 3.419508563221334
$ws.Range("R6").Value = 30.775577068992
$ws.Range("S6").Value = 0.005835310823747441
$ws.Range("T6").Value = 0.005835310823747441

$ws.Range("G7").Value = 0.100996
$ws.Range("H7").Value = 0.302988
$ws.Range("I7").Value = 0.009287810103293732
$ws.Range("J7").Value = 0.009287810103293733
$ws.Range("M7").Value = 18.53974466666667
$ws.Range("N7").Value = 55.619234
$ws.Range("O7").Value = 0.3440288735868225
$ws.Range("P7").Value = 0.3440288735868225
$ws.Range("Q7").Value = 1.872440052354666
$ws.Range("R7").Value = 16.851960471192
$ws.Range("S7").Value = 0.003195274847924453
$ws.Range("T7").Value = 0.003195274847924453

$ws.Range("G8").Value = 0.6471736666666666
$ws.Range("H8").Value = 1.941521
$ws.Range("I8").Value = 0.05951548694851595
$ws.Range("J8").Value = 0.05951548694851596
$ws.Range("M8").Value = 1.492477333333333
$ws.Range("N8").Value = 4.477432
$ws.Range("O8").Value = 0.02769484181536182
$ws.Range("P8").Value = 0.02769484181536182
$ws.Range("Q8").Value = 0.9658920282302222
$ws.Range("R8").Value = 8.693028254072001
$ws.Range("S8").Value = 0.001648271996603381
$ws.Range("T8").Value = 0.001648271996603381

$ws.Range("G9").Value = 0.6471736666666666
$ws.Range("H9").Value = 1.941521
$ws.Range("I9").Value = 0.05951548694851595
$ws.Range("J9").Value = 0.05951548694851596
$ws.Range("M9").Value = 33.85786133333334
$ws.Range("N9").Value = 101.573584
$ws.Range("O9").Value = 0.6282762845978157
$ws.Range("P9").Value = 0.6282762845978156
$ws.Range("Q9").Value = 21.91191626458489
$ws.Range("R9").Value = 197.207246381264
$ws.Range("S9").Value = 0.03739216901604339
$ws.Range("T9").Value = 0.03739216901604339

$ws.Range("G10").Value = 0.6471736666666666
$ws.Range("H10").Value = 1.941521
$ws.Range("I10").Value = 0.05951548694851595
$ws.Range("J10").Value = 0.05951548694851596
$ws.Range("M10").Value = 18.53974466666667
$ws.Range("N10").Value = 55.619234
$ws.Range("O10").Value = 0.3440288735868225
$ws.Range("P10").Value = 0.3440288735868225
$ws.Range("Q10").Value = 11.99843453499044
$ws.Range("R10").Value = 107.985910814914
$ws.Range("S10").Value = 0.02047504593586918
$ws.Range("T10").Value = 0.02047504593586918

